$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the tiny float drift on the previous row's timestamp (recorded at
# save time rather than at the moment the prior operation ran).
$ws.Range("D6").Value = 45735.95928607639

# New row for the "Raiz" (square root) operation.
$ws.Range("A7").Value = "Marvin"
$ws.Range("B7").Value = "125 Raíz 2"
$ws.Range("C7").Value = 11.18033988749895
$ws.Range("D7").Value = 45735.96277760838
$ws.Range("D7").NumberFormat = $ws.Range("D6").NumberFormat
